# Fixing network data cleaning scripts
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header columns (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Title-case the connector words ("de", "del", "la", "los") within
#    state/municipality names.
$ws.Range("B8").Value  = "Amatenango De La Frontera"
$ws.Range("B10").Value = "Comitán De Domínguez"
$ws.Range("B16").Value = "Salto De Agua"
$ws.Range("A28").Value = "Ciudad De México"
$ws.Range("A42").Value = "Estado De México"
$ws.Range("B42").Value = "Almoloya De Alquisiras"
$ws.Range("B44").Value = "Ecatepec De Morelos"
$ws.Range("B50").Value = "Tlalnepantla De Baz"
$ws.Range("B52").Value = "San Miguel De Allende"
$ws.Range("B55").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B57").Value = "Acapulco De Juárez"
$ws.Range("B61").Value = "Atoyac De Álvarez"
$ws.Range("B62").Value = "Chilapa De Álvarez"
$ws.Range("B65").Value = "Coyuca De Benítez"
$ws.Range("B69").Value = "Tlapa De Comonfort"
$ws.Range("B73").Value = "Pachuca De Soto"
$ws.Range("B78").Value = "Lagos De Moreno"
$ws.Range("B97").Value = "San Nicolás De Los Garza"
$ws.Range("B100").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B101").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B103").Value = "Ocotlán De Morelos"
$ws.Range("B112").Value = "Teotitlán Del Valle"
$ws.Range("B115").Value = "Cuayuca De Andrade"
$ws.Range("B117").Value = "Ixcamilpa De Guerrero"
$ws.Range("B129").Value = "Tepanco De López"
$ws.Range("B133").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B141").Value = "Landa De Matamoros"
$ws.Range("B159").Value = "Chinampa De Gorostiza"
$ws.Range("B161").Value = "Hueyapan De Ocampo"

# 3) Remove trailing metadata/footer rows 177-181 (row 176 is already blank,
#    leaving the sheet data ending at row 175).
$ws.Range("A177:A181").EntireRow.Delete()
